# Updates the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# with refreshed values, matching the GitHub Actions scheduled data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold display strings (e.g. "1.007", "26.988.68",
# "  +2.15%  ") that must stay text, not get auto-coerced to numbers/dates by
# Excel's smart input. Forcing NumberFormat to "@" (Text) before assigning the
# value keeps it a string, then resetting the style back to "Normal" avoids
# leaving a stray text-format style applied to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.988.68"
Set-TextValue $ws.Range("E2") "  +2.15%  "
Set-TextValue $ws.Range("D3") "1.849.52"
Set-TextValue $ws.Range("E3") "  +2.40%  "
Set-TextValue $ws.Range("D4") "1.007"
Set-TextValue $ws.Range("E4") "  -0.22%  "
Set-TextValue $ws.Range("D5") "310.72"
Set-TextValue $ws.Range("E5") "  +1.24%  "
Set-TextValue $ws.Range("E6") "  -0.05%  "
Set-TextValue $ws.Range("E7") "  +3.40%  "
Set-TextValue $ws.Range("D8") "0.3634"
Set-TextValue $ws.Range("E8") "  +0.98%  "
Set-TextValue $ws.Range("D9") "0.07179"
Set-TextValue $ws.Range("E9") "  +1.53%  "
Set-TextValue $ws.Range("E10") "  +4.32%  "
Set-TextValue $ws.Range("D11") "19.63"
Set-TextValue $ws.Range("E11") "  +1.02%  "
Set-TextValue $ws.Range("D12") "0.07683"
Set-TextValue $ws.Range("E12") "  -1.70%  "
Set-TextValue $ws.Range("D13") "1.861.37"
Set-TextValue $ws.Range("E13") "  +3.05%  "
Set-TextValue $ws.Range("D14") "5.298"
Set-TextValue $ws.Range("E14") "  +0.15%  "
Set-TextValue $ws.Range("D15") "6.417"
Set-TextValue $ws.Range("E15") "  +1.47%  "
Set-TextValue $ws.Range("D16") "88.30"
Set-TextValue $ws.Range("E16") "  +3.50%  "
Set-TextValue $ws.Range("D17") "1.009"
Set-TextValue $ws.Range("E17") "  -0.14%  "
Set-TextValue $ws.Range("D18") "0.000008608"
Set-TextValue $ws.Range("E18") "  +1.24%  "
Set-TextValue $ws.Range("D19") "1.007"
Set-TextValue $ws.Range("E19") "  +0.08%  "
Set-TextValue $ws.Range("D20") "27.023.38"
Set-TextValue $ws.Range("E20") "  +2.13%  "
Set-TextValue $ws.Range("E21") "  +1.60%  "
Set-TextValue $ws.Range("D22") "5.039"
Set-TextValue $ws.Range("E22") "  +1.40%  "
Set-TextValue $ws.Range("E23") "  +1.27%  "
Set-TextValue $ws.Range("D24") "1.939"
Set-TextValue $ws.Range("D25") "152.66"
Set-TextValue $ws.Range("E25") "  -0.10%  "
Set-TextValue $ws.Range("D26") "18.13"
Set-TextValue $ws.Range("E26") "  +1.89%  "
Set-TextValue $ws.Range("D27") "2.044"
Set-TextValue $ws.Range("E27") "  -1.62%  "
Set-TextValue $ws.Range("E28") "  +1.75%  "
Set-TextValue $ws.Range("D29") "4.935"
Set-TextValue $ws.Range("E29") "  +1.65%  "
Set-TextValue $ws.Range("D30") "0.08864"
Set-TextValue $ws.Range("E30") "  +1.92%  "
Set-TextValue $ws.Range("E31") "  +2.06%  "
Set-TextValue $ws.Range("D32") "2.840"
Set-TextValue $ws.Range("E32") "  +0.94%  "
Set-TextValue $ws.Range("D33") "1.181"
Set-TextValue $ws.Range("E33") "  +6.84%  "
Set-TextValue $ws.Range("D34") "0.7458"
Set-TextValue $ws.Range("E34") "  +2.82%  "
Set-TextValue $ws.Range("E35") "  +0.43%  "
Set-TextValue $ws.Range("E36") "  +0.95%  "
Set-TextValue $ws.Range("D37") "2.989"
Set-TextValue $ws.Range("E37") "  +2.47%  "
Set-TextValue $ws.Range("D38") "0.01944"
Set-TextValue $ws.Range("E38") "  +0.34%  "
Set-TextValue $ws.Range("D39") "0.05173"
Set-TextValue $ws.Range("E39") "  +1.04%  "
Set-TextValue $ws.Range("D40") "0.5161"
Set-TextValue $ws.Range("E40") "  +1.54%  "
Set-TextValue $ws.Range("D41") "6.903"
Set-TextValue $ws.Range("E41") "  +1.79%  "
Set-TextValue $ws.Range("E42") "  +0.04%  "
Set-TextValue $ws.Range("D43") "8.181"
Set-TextValue $ws.Range("E43") "  +1.94%  "
Set-TextValue $ws.Range("D44") "10.54"
Set-TextValue $ws.Range("E44") "  +4.96%  "
Set-TextValue $ws.Range("D45") "0.4721"
Set-TextValue $ws.Range("E45") "  +1.10%  "
Set-TextValue $ws.Range("D46") "1.008"
Set-TextValue $ws.Range("E46") "  +0.00%  "
Set-TextValue $ws.Range("D47") "100.35"
Set-TextValue $ws.Range("E47") "  -0.13%  "
Set-TextValue $ws.Range("D48") "1.605"
Set-TextValue $ws.Range("E48") "  +1.78%  "
Set-TextValue $ws.Range("D49") "0.06053"
Set-TextValue $ws.Range("E49") "  +1.10%  "
Set-TextValue $ws.Range("D50") "64.43"
Set-TextValue $ws.Range("E50") "  +0.99%  "
Set-TextValue $ws.Range("E51") "  -0.01%  "
